# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" wherever it is used
#    (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Latest Handoff Datetime"-ish status columns:
#    Overview columns E & F, zh-cn column C, de-de column C.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update status text --------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow columns -------------------------------------------------------
# The engine snaps ColumnWidth to an internal 1/6-character grid
# (stored_width = round6(ColumnWidth + 5/6)); feed it the inverse of the
# target stored width (13.4101845877511 chars) so the saved file lands on
# the closest achievable grid point.
$targetWidth = 13.4101845877511
$newColumnWidth = $targetWidth - (5.0/6.0)

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
